$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.933367848396301
$ws.Range("B1").Value = 1.994569778442383
$ws.Range("C1").Value = 2.031540393829346
$ws.Range("D1").Value = 2.565282821655273
$ws.Range("E1").Value = 3.325434684753418
